$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("I4").Value = -0.1931047144327512
$ws.Range("J4").Value = 0.4780943459631815
$ws.Range("K4").Value = 0.4160877221743743
$ws.Range("L4").Value = 2.714819154945082
